$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1016.55554
$ws.Range("I40").Value = 920
$ws.Range("J40").Value = 1137.25
$ws.Range("K40").Value = 920
$ws.Range("L40").Value = 1137.25
$ws.Range("M40").Value = -745
$ws.Range("N40").Value = -1487.25
$ws.Range("H74").Value = 3616.3635
$ws.Range("J74").Value = 3578
$ws.Range("L74").Value = 3578
$ws.Range("N74").Value = -5450
$ws.Range("H77").Value = 3616.3635
$ws.Range("J77").Value = 3578
$ws.Range("L77").Value = 17890
$ws.Range("N77").Value = -27250
$ws.Range("H112").Value = 1356.1224
$ws.Range("J112").Value = 1365.9794
$ws.Range("L112").Value = 4097.9382
$ws.Range("N112").Value = -6313.9382
$ws.Range("H133").Value = 25981.54
$ws.Range("J133").Value = 25981.54
$ws.Range("L133").Value = 25981.54
$ws.Range("N133").Value = -36101.54
$ws.Range("H134").Value = 30556
$ws.Range("J134").Value = 30556
$ws.Range("L134").Value = 30556
$ws.Range("N134").Value = -40696
$ws.Range("H136").Value = 29150
$ws.Range("J136").Value = 29150
$ws.Range("L136").Value = 29150
$ws.Range("N136").Value = -39350
$ws.Range("H137").Value = 2529.8125
$ws.Range("I137").Value = 2657.9429
$ws.Range("J137").Value = 2184.8462
$ws.Range("K137").Value = 7973.8287
$ws.Range("L137").Value = 6554.5386
$ws.Range("M137").Value = -5423.8287
$ws.Range("N137").Value = -11654.5386
$ws.Range("H138").Value = 4524.057
$ws.Range("I138").Value = 1502.6666
$ws.Range("J138").Value = 5818.939
$ws.Range("K138").Value = 4507.9998
$ws.Range("L138").Value = 17456.817
$ws.Range("M138").Value = 632.0002000000004
$ws.Range("N138").Value = -27736.817
$ws.Range("H139").Value = 30000
$ws.Range("J139").Value = 30000
$ws.Range("L139").Value = 30000
$ws.Range("N139").Value = -40280
$ws.Range("H141").Value = 320794.62
$ws.Range("I141").Value = 1194.7931
$ws.Range("J141").Value = 1644851
$ws.Range("K141").Value = 3584.379300000001
$ws.Range("L141").Value = 4934553
$ws.Range("M141").Value = 1595.620699999999
$ws.Range("N141").Value = -4944913

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H6").Value = 215001.5
$ws.Range("I6").Value = 525000
$ws.Range("J6").Value = 60002.25
$ws.Range("K6").Value = 525000
$ws.Range("L6").Value = 60002.25
$ws.Range("M6").Value = -524827
$ws.Range("N6").Value = -60348.25
$ws.Range("H32").Value = 1490.48
$ws.Range("I32").Value = 1226.1555
$ws.Range("K32").Value = 1226.1555
$ws.Range("M32").Value = -939.1555000000001
$ws.Range("H74").Value = 939.61536
$ws.Range("I74").Value = 846.1111
$ws.Range("K74").Value = 846.1111
$ws.Range("M74").Value = 27.88890000000004
$ws.Range("H77").Value = 939.61536
$ws.Range("I77").Value = 846.1111
$ws.Range("K77").Value = 4230.555499999999
$ws.Range("M77").Value = 137.4445000000005
$ws.Range("H97").Value = 787.86664
$ws.Range("I97").Value = 770.61536
$ws.Range("K97").Value = 770.61536
$ws.Range("M97").Value = -274.61536
$ws.Range("H132").Value = 19232998
$ws.Range("I132").Value = 25642514
$ws.Range("K132").Value = 76927542
$ws.Range("M132").Value = -76925012
$ws.Range("H139").Value = 26123.777
$ws.Range("J139").Value = 26123.777
$ws.Range("L139").Value = 26123.777
$ws.Range("N139").Value = -36403.777

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5335.579
$ws.Range("I20").Value = 5511.9287
$ws.Range("J20").Value = 4841.8
$ws.Range("K20").Value = 5511.9287
$ws.Range("L20").Value = 4841.8
$ws.Range("M20").Value = -5264.9287
$ws.Range("N20").Value = -5335.8
$ws.Range("H82").Value = 20640.785
$ws.Range("I82").Value = 2708.8572
$ws.Range("J82").Value = 38572.715
$ws.Range("K82").Value = 2708.8572
$ws.Range("L82").Value = 38572.715
$ws.Range("M82").Value = -2325.8572
$ws.Range("N82").Value = -39338.715
$ws.Range("H85").Value = 20640.785
$ws.Range("I85").Value = 2708.8572
$ws.Range("J85").Value = 38572.715
$ws.Range("K85").Value = 2708.8572
$ws.Range("L85").Value = 38572.715
$ws.Range("M85").Value = -1382.8572
$ws.Range("N85").Value = -41224.715
$ws.Range("H94").Value = 734.4138
$ws.Range("I94").Value = 658.2353000000001
$ws.Range("J94").Value = 842.3333
$ws.Range("K94").Value = 658.2353000000001
$ws.Range("L94").Value = 842.3333
$ws.Range("M94").Value = -207.2353000000001
$ws.Range("N94").Value = -1744.3333
$ws.Range("H134").Value = 4116.2856
$ws.Range("I134").Value = 3136
$ws.Range("J134").Value = 5880.8
$ws.Range("K134").Value = 9408
$ws.Range("L134").Value = 17642.4
$ws.Range("M134").Value = -6873
$ws.Range("N134").Value = -22712.4
$ws.Range("H141").Value = 37045.453
$ws.Range("J141").Value = 29642.857
$ws.Range("L141").Value = 29642.857
$ws.Range("N141").Value = -40002.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 70013
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H31").Value = 2633.3726
$ws.Range("I31").Value = 1564.2963
$ws.Range("J31").Value = 3836.0833
$ws.Range("K31").Value = 1564.2963
$ws.Range("L31").Value = 3836.0833
$ws.Range("M31").Value = -1269.2963
$ws.Range("N31").Value = -4426.0833
$ws.Range("H34").Value = 2633.3726
$ws.Range("I34").Value = 1564.2963
$ws.Range("J34").Value = 3836.0833
$ws.Range("K34").Value = 1564.2963
$ws.Range("L34").Value = 3836.0833
$ws.Range("M34").Value = -1362.2963
$ws.Range("N34").Value = -4240.0833
$ws.Range("H50").Value = 16564
$ws.Range("J50").Value = 16564
$ws.Range("L50").Value = 16564
$ws.Range("N50").Value = -17814
$ws.Range("H105").Value = 3460.9
$ws.Range("I105").Value = 3451.2856
$ws.Range("J105").Value = 3483.3333
$ws.Range("K105").Value = 3451.2856
$ws.Range("L105").Value = 3483.3333
$ws.Range("M105").Value = -1704.2856
$ws.Range("N105").Value = -6977.3333
$ws.Range("H127").Value = 32993.332
$ws.Range("J127").Value = 32993.332
$ws.Range("L127").Value = 32993.332
$ws.Range("N127").Value = -42913.332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6129500
$ws.Range("I11").Value = 6600000
$ws.Range("J11").Value = 5345334
$ws.Range("K11").Value = 6600000
$ws.Range("L11").Value = 5345334
$ws.Range("M11").Value = -6599861
$ws.Range("N11").Value = -5345612
$ws.Range("H18").Value = 200575.72
$ws.Range("J18").Value = 67338.336
$ws.Range("L18").Value = 67338.336
$ws.Range("N18").Value = -67924.336
$ws.Range("H97").Value = 1932
$ws.Range("I97").Value = 1156.9333
$ws.Range("K97").Value = 1156.9333
$ws.Range("M97").Value = -660.9332999999999
$ws.Range("H132").Value = 4867.222
$ws.Range("I132").Value = 6922.4
$ws.Range("J132").Value = 4076.7693
$ws.Range("K132").Value = 20767.2
$ws.Range("L132").Value = 12230.3079
$ws.Range("M132").Value = -18237.2
$ws.Range("N132").Value = -17290.3079
$ws.Range("H139").Value = 29975
$ws.Range("J139").Value = 29975
$ws.Range("L139").Value = 29975
$ws.Range("N139").Value = -40255

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1405.25
$ws.Range("I9").Value = 280
$ws.Range("J9").Value = 2530.5
$ws.Range("K9").Value = 280
$ws.Range("L9").Value = 2530.5
$ws.Range("M9").Value = -56
$ws.Range("N9").Value = -2978.5
$ws.Range("H40").Value = 2626
$ws.Range("I40").Value = 1004
$ws.Range("J40").Value = 3166.6667
$ws.Range("K40").Value = 1004
$ws.Range("L40").Value = 3166.6667
$ws.Range("M40").Value = -868
$ws.Range("N40").Value = -3438.6667
$ws.Range("H51").Value = 25084
$ws.Range("J51").Value = 25084
$ws.Range("L51").Value = 25084
$ws.Range("N51").Value = -26040
$ws.Range("H132").Value = 3252.5
$ws.Range("I132").Value = 1804.1875
$ws.Range("J132").Value = 5183.5835
$ws.Range("K132").Value = 5412.5625
$ws.Range("L132").Value = 15550.7505
$ws.Range("M132").Value = -2882.5625
$ws.Range("N132").Value = -20610.7505
$ws.Range("H140").Value = 29551.666
$ws.Range("J140").Value = 29551.666
$ws.Range("L140").Value = 29551.666
$ws.Range("N140").Value = -39911.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10486.099
$ws.Range("I132").Value = 2116.451
$ws.Range("K132").Value = 6349.353
$ws.Range("M132").Value = -3819.353
$ws.Range("H136").Value = 1010.95557
$ws.Range("I136").Value = 630.40625
$ws.Range("J136").Value = 1947.6923
$ws.Range("K136").Value = 1891.21875
$ws.Range("L136").Value = 5843.0769
$ws.Range("M136").Value = 658.78125
$ws.Range("N136").Value = -10943.0769
$ws.Range("H138").Value = 30000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 30000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 30000
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -40280
